$wb = $excel.ActiveWorkbook

# --- Proximity sheet: add rows 8-9 ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(8, 1).Formula = "'2026-01-28"
$ws.Cells.Item(8, 2).Value = "17:41:16"
$ws.Cells.Item(8, 3).Value = "17:00"
$ws.Cells.Item(8, 4).Value = "Living Room Main Door"
$ws.Cells.Item(8, 5).Value = "ENTER"
$ws.Cells.Item(8, 6).Value = "User ENTERED Living Room Main Door"

$ws.Cells.Item(9, 1).Formula = "'2026-01-28"
$ws.Cells.Item(9, 2).Value = "17:41:18"
$ws.Cells.Item(9, 3).Value = "17:00"
$ws.Cells.Item(9, 4).Value = "Living Room Main Door"
$ws.Cells.Item(9, 5).Value = "EXIT"
$ws.Cells.Item(9, 6).Value = "User EXITED Living Room Main Door"

# --- mmWave sheet: add rows 97-148 ---
$ws = $wb.Worksheets.Item("mmWave")
$ws.Cells.Item(97, 1).Formula = "'2026-01-28"
$ws.Cells.Item(97, 2).Value = "17:40:59"
$ws.Cells.Item(97, 3).Value = "17:00"
$ws.Cells.Item(97, 4).Value = "Bedroom"
$ws.Cells.Item(97, 5).Value = "In Bed | HR=0 | BR=0"
$ws.Cells.Item(97, 6).Value = "Occupied"

$ws.Cells.Item(98, 1).Formula = "'2026-01-28"
$ws.Cells.Item(98, 2).Value = "17:40:59"
$ws.Cells.Item(98, 3).Value = "17:00"
$ws.Cells.Item(98, 4).Value = "Bedroom"
$ws.Cells.Item(98, 5).Value = "In Bed | HR=106 | BR=58"
$ws.Cells.Item(98, 6).Value = "Occupied"

$ws.Cells.Item(99, 1).Formula = "'2026-01-28"
$ws.Cells.Item(99, 2).Value = "17:41:00"
$ws.Cells.Item(99, 3).Value = "17:00"
$ws.Cells.Item(99, 4).Value = "Bedroom"
$ws.Cells.Item(99, 5).Value = "In Bed | HR=90 | BR=42"
$ws.Cells.Item(99, 6).Value = "Occupied"

$ws.Cells.Item(100, 1).Formula = "'2026-01-28"
$ws.Cells.Item(100, 2).Value = "17:41:00"
$ws.Cells.Item(100, 3).Value = "17:00"
$ws.Cells.Item(100, 4).Value = "Bedroom"
$ws.Cells.Item(100, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(100, 6).Value = "Occupied"

$ws.Cells.Item(101, 1).Formula = "'2026-01-28"
$ws.Cells.Item(101, 2).Value = "17:41:00"
$ws.Cells.Item(101, 3).Value = "17:00"
$ws.Cells.Item(101, 4).Value = "Bedroom"
$ws.Cells.Item(101, 5).Value = "In Bed | HR=53 | BR=5"
$ws.Cells.Item(101, 6).Value = "Occupied"

$ws.Cells.Item(102, 1).Formula = "'2026-01-28"
$ws.Cells.Item(102, 2).Value = "17:41:00"
$ws.Cells.Item(102, 3).Value = "17:00"
$ws.Cells.Item(102, 4).Value = "Bedroom"
$ws.Cells.Item(102, 5).Value = "In Bed | HR=58 | BR=10"
$ws.Cells.Item(102, 6).Value = "Occupied"

$ws.Cells.Item(103, 1).Formula = "'2026-01-28"
$ws.Cells.Item(103, 2).Value = "17:41:00"
$ws.Cells.Item(103, 3).Value = "17:00"
$ws.Cells.Item(103, 4).Value = "Bedroom"
$ws.Cells.Item(103, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(103, 6).Value = "Occupied"

$ws.Cells.Item(104, 1).Formula = "'2026-01-28"
$ws.Cells.Item(104, 2).Value = "17:41:05"
$ws.Cells.Item(104, 3).Value = "17:00"
$ws.Cells.Item(104, 4).Value = "Bedroom"
$ws.Cells.Item(104, 5).Value = "In Bed | HR=59 | BR=11"
$ws.Cells.Item(104, 6).Value = "Occupied"

$ws.Cells.Item(105, 1).Formula = "'2026-01-28"
$ws.Cells.Item(105, 2).Value = "17:41:05"
$ws.Cells.Item(105, 3).Value = "17:00"
$ws.Cells.Item(105, 4).Value = "Bedroom"
$ws.Cells.Item(105, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(105, 6).Value = "Occupied"

$ws.Cells.Item(106, 1).Formula = "'2026-01-28"
$ws.Cells.Item(106, 2).Value = "17:41:07"
$ws.Cells.Item(106, 3).Value = "17:00"
$ws.Cells.Item(106, 4).Value = "Bedroom"
$ws.Cells.Item(106, 5).Value = "In Bed | HR=59 | BR=11"
$ws.Cells.Item(106, 6).Value = "Occupied"

$ws.Cells.Item(107, 1).Formula = "'2026-01-28"
$ws.Cells.Item(107, 2).Value = "17:41:08"
$ws.Cells.Item(107, 3).Value = "17:00"
$ws.Cells.Item(107, 4).Value = "Bedroom"
$ws.Cells.Item(107, 5).Value = "In Bed | HR=101 | BR=53"
$ws.Cells.Item(107, 6).Value = "Occupied"

$ws.Cells.Item(108, 1).Formula = "'2026-01-28"
$ws.Cells.Item(108, 2).Value = "17:41:09"
$ws.Cells.Item(108, 3).Value = "17:00"
$ws.Cells.Item(108, 4).Value = "Bedroom"
$ws.Cells.Item(108, 5).Value = "In Bed | HR=87 | BR=39"
$ws.Cells.Item(108, 6).Value = "Occupied"

$ws.Cells.Item(109, 1).Formula = "'2026-01-28"
$ws.Cells.Item(109, 2).Value = "17:41:10"
$ws.Cells.Item(109, 3).Value = "17:00"
$ws.Cells.Item(109, 4).Value = "Bedroom"
$ws.Cells.Item(109, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(109, 6).Value = "Occupied"

$ws.Cells.Item(110, 1).Formula = "'2026-01-28"
$ws.Cells.Item(110, 2).Value = "17:41:11"
$ws.Cells.Item(110, 3).Value = "17:00"
$ws.Cells.Item(110, 4).Value = "Bedroom"
$ws.Cells.Item(110, 5).Value = "In Bed | HR=52 | BR=4"
$ws.Cells.Item(110, 6).Value = "Occupied"

$ws.Cells.Item(111, 1).Formula = "'2026-01-28"
$ws.Cells.Item(111, 2).Value = "17:41:11"
$ws.Cells.Item(111, 3).Value = "17:00"
$ws.Cells.Item(111, 4).Value = "Bedroom"
$ws.Cells.Item(111, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(111, 6).Value = "Occupied"

$ws.Cells.Item(112, 1).Formula = "'2026-01-28"
$ws.Cells.Item(112, 2).Value = "17:41:13"
$ws.Cells.Item(112, 3).Value = "17:00"
$ws.Cells.Item(112, 4).Value = "Bedroom"
$ws.Cells.Item(112, 5).Value = "In Bed | HR=53 | BR=5"
$ws.Cells.Item(112, 6).Value = "Occupied"

$ws.Cells.Item(113, 1).Formula = "'2026-01-28"
$ws.Cells.Item(113, 2).Value = "17:41:14"
$ws.Cells.Item(113, 3).Value = "17:00"
$ws.Cells.Item(113, 4).Value = "Bedroom"
$ws.Cells.Item(113, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(113, 6).Value = "Occupied"

$ws.Cells.Item(114, 1).Formula = "'2026-01-28"
$ws.Cells.Item(114, 2).Value = "17:41:16"
$ws.Cells.Item(114, 3).Value = "17:00"
$ws.Cells.Item(114, 4).Value = "Bedroom"
$ws.Cells.Item(114, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(114, 6).Value = "Occupied"

$ws.Cells.Item(115, 1).Formula = "'2026-01-28"
$ws.Cells.Item(115, 2).Value = "17:41:17"
$ws.Cells.Item(115, 3).Value = "17:00"
$ws.Cells.Item(115, 4).Value = "Bedroom"
$ws.Cells.Item(115, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(115, 6).Value = "Occupied"

$ws.Cells.Item(116, 1).Formula = "'2026-01-28"
$ws.Cells.Item(116, 2).Value = "17:41:20"
$ws.Cells.Item(116, 3).Value = "17:00"
$ws.Cells.Item(116, 4).Value = "Bedroom"
$ws.Cells.Item(116, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(116, 6).Value = "Occupied"

$ws.Cells.Item(117, 1).Formula = "'2026-01-28"
$ws.Cells.Item(117, 2).Value = "17:41:22"
$ws.Cells.Item(117, 3).Value = "17:00"
$ws.Cells.Item(117, 4).Value = "Bedroom"
$ws.Cells.Item(117, 5).Value = "In Bed | HR=111 | BR=63"
$ws.Cells.Item(117, 6).Value = "Occupied"

$ws.Cells.Item(118, 1).Formula = "'2026-01-28"
$ws.Cells.Item(118, 2).Value = "17:41:23"
$ws.Cells.Item(118, 3).Value = "17:00"
$ws.Cells.Item(118, 4).Value = "Bedroom"
$ws.Cells.Item(118, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(118, 6).Value = "Occupied"

$ws.Cells.Item(119, 1).Formula = "'2026-01-28"
$ws.Cells.Item(119, 2).Value = "17:41:23"
$ws.Cells.Item(119, 3).Value = "17:00"
$ws.Cells.Item(119, 4).Value = "Bedroom"
$ws.Cells.Item(119, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(119, 6).Value = "Occupied"

$ws.Cells.Item(120, 1).Formula = "'2026-01-28"
$ws.Cells.Item(120, 2).Value = "17:41:25"
$ws.Cells.Item(120, 3).Value = "17:00"
$ws.Cells.Item(120, 4).Value = "Bedroom"
$ws.Cells.Item(120, 5).Value = "In Bed | HR=92 | BR=44"
$ws.Cells.Item(120, 6).Value = "Occupied"

$ws.Cells.Item(121, 1).Formula = "'2026-01-28"
$ws.Cells.Item(121, 2).Value = "17:41:26"
$ws.Cells.Item(121, 3).Value = "17:00"
$ws.Cells.Item(121, 4).Value = "Bedroom"
$ws.Cells.Item(121, 5).Value = "In Bed | HR=74 | BR=26"
$ws.Cells.Item(121, 6).Value = "Occupied"

$ws.Cells.Item(122, 1).Formula = "'2026-01-28"
$ws.Cells.Item(122, 2).Value = "17:41:26"
$ws.Cells.Item(122, 3).Value = "17:00"
$ws.Cells.Item(122, 4).Value = "Bedroom"
$ws.Cells.Item(122, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(122, 6).Value = "Occupied"

$ws.Cells.Item(123, 1).Formula = "'2026-01-28"
$ws.Cells.Item(123, 2).Value = "17:41:28"
$ws.Cells.Item(123, 3).Value = "17:00"
$ws.Cells.Item(123, 4).Value = "Bedroom"
$ws.Cells.Item(123, 5).Value = "In Bed | HR=83 | BR=35"
$ws.Cells.Item(123, 6).Value = "Occupied"

$ws.Cells.Item(124, 1).Formula = "'2026-01-28"
$ws.Cells.Item(124, 2).Value = "17:41:29"
$ws.Cells.Item(124, 3).Value = "17:00"
$ws.Cells.Item(124, 4).Value = "Bedroom"
$ws.Cells.Item(124, 5).Value = "In Bed | HR=92 | BR=44"
$ws.Cells.Item(124, 6).Value = "Occupied"

$ws.Cells.Item(125, 1).Formula = "'2026-01-28"
$ws.Cells.Item(125, 2).Value = "17:41:29"
$ws.Cells.Item(125, 3).Value = "17:00"
$ws.Cells.Item(125, 4).Value = "Bedroom"
$ws.Cells.Item(125, 5).Value = "In Bed | HR=56 | BR=8"
$ws.Cells.Item(125, 6).Value = "Occupied"

$ws.Cells.Item(126, 1).Formula = "'2026-01-28"
$ws.Cells.Item(126, 2).Value = "17:41:31"
$ws.Cells.Item(126, 3).Value = "17:00"
$ws.Cells.Item(126, 4).Value = "Bedroom"
$ws.Cells.Item(126, 5).Value = "In Bed | HR=127 | BR=79"
$ws.Cells.Item(126, 6).Value = "Occupied"

$ws.Cells.Item(127, 1).Formula = "'2026-01-28"
$ws.Cells.Item(127, 2).Value = "17:41:32"
$ws.Cells.Item(127, 3).Value = "17:00"
$ws.Cells.Item(127, 4).Value = "Bedroom"
$ws.Cells.Item(127, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(127, 6).Value = "Occupied"

$ws.Cells.Item(128, 1).Formula = "'2026-01-28"
$ws.Cells.Item(128, 2).Value = "17:41:34"
$ws.Cells.Item(128, 3).Value = "17:00"
$ws.Cells.Item(128, 4).Value = "Bedroom"
$ws.Cells.Item(128, 5).Value = "In Bed | HR=55 | BR=7"
$ws.Cells.Item(128, 6).Value = "Occupied"

$ws.Cells.Item(129, 1).Formula = "'2026-01-28"
$ws.Cells.Item(129, 2).Value = "17:41:35"
$ws.Cells.Item(129, 3).Value = "17:00"
$ws.Cells.Item(129, 4).Value = "Bedroom"
$ws.Cells.Item(129, 5).Value = "In Bed | HR=113 | BR=65"
$ws.Cells.Item(129, 6).Value = "Occupied"

$ws.Cells.Item(130, 1).Formula = "'2026-01-28"
$ws.Cells.Item(130, 2).Value = "17:41:35"
$ws.Cells.Item(130, 3).Value = "17:00"
$ws.Cells.Item(130, 4).Value = "Bedroom"
$ws.Cells.Item(130, 5).Value = "In Bed | HR=57 | BR=9"
$ws.Cells.Item(130, 6).Value = "Occupied"

$ws.Cells.Item(131, 1).Formula = "'2026-01-28"
$ws.Cells.Item(131, 2).Value = "17:41:37"
$ws.Cells.Item(131, 3).Value = "17:00"
$ws.Cells.Item(131, 4).Value = "Bedroom"
$ws.Cells.Item(131, 5).Value = "In Bed | HR=91 | BR=43"
$ws.Cells.Item(131, 6).Value = "Occupied"

$ws.Cells.Item(132, 1).Formula = "'2026-01-28"
$ws.Cells.Item(132, 2).Value = "17:41:37"
$ws.Cells.Item(132, 3).Value = "17:00"
$ws.Cells.Item(132, 4).Value = "Bedroom"
$ws.Cells.Item(132, 5).Value = "In Bed | HR=55 | BR=7"
$ws.Cells.Item(132, 6).Value = "Occupied"

$ws.Cells.Item(133, 1).Formula = "'2026-01-28"
$ws.Cells.Item(133, 2).Value = "17:41:38"
$ws.Cells.Item(133, 3).Value = "17:00"
$ws.Cells.Item(133, 4).Value = "Bedroom"
$ws.Cells.Item(133, 5).Value = "In Bed | HR=53 | BR=5"
$ws.Cells.Item(133, 6).Value = "Occupied"

$ws.Cells.Item(134, 1).Formula = "'2026-01-28"
$ws.Cells.Item(134, 2).Value = "17:41:40"
$ws.Cells.Item(134, 3).Value = "17:00"
$ws.Cells.Item(134, 4).Value = "Bedroom"
$ws.Cells.Item(134, 5).Value = "In Bed | HR=62 | BR=14"
$ws.Cells.Item(134, 6).Value = "Occupied"

$ws.Cells.Item(135, 1).Formula = "'2026-01-28"
$ws.Cells.Item(135, 2).Value = "17:41:40"
$ws.Cells.Item(135, 3).Value = "17:00"
$ws.Cells.Item(135, 4).Value = "Bedroom"
$ws.Cells.Item(135, 5).Value = "In Bed | HR=53 | BR=5"
$ws.Cells.Item(135, 6).Value = "Occupied"

$ws.Cells.Item(136, 1).Formula = "'2026-01-28"
$ws.Cells.Item(136, 2).Value = "17:41:41"
$ws.Cells.Item(136, 3).Value = "17:00"
$ws.Cells.Item(136, 4).Value = "Bedroom"
$ws.Cells.Item(136, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(136, 6).Value = "Occupied"

$ws.Cells.Item(137, 1).Formula = "'2026-01-28"
$ws.Cells.Item(137, 2).Value = "17:41:43"
$ws.Cells.Item(137, 3).Value = "17:00"
$ws.Cells.Item(137, 4).Value = "Bedroom"
$ws.Cells.Item(137, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(137, 6).Value = "Occupied"

$ws.Cells.Item(138, 1).Formula = "'2026-01-28"
$ws.Cells.Item(138, 2).Value = "17:41:47"
$ws.Cells.Item(138, 3).Value = "17:00"
$ws.Cells.Item(138, 4).Value = "Bedroom"
$ws.Cells.Item(138, 5).Value = "In Bed | HR=53 | BR=5"
$ws.Cells.Item(138, 6).Value = "Occupied"

$ws.Cells.Item(139, 1).Formula = "'2026-01-28"
$ws.Cells.Item(139, 2).Value = "17:41:49"
$ws.Cells.Item(139, 3).Value = "17:00"
$ws.Cells.Item(139, 4).Value = "Bedroom"
$ws.Cells.Item(139, 5).Value = "In Bed | HR=60 | BR=12"
$ws.Cells.Item(139, 6).Value = "Occupied"

$ws.Cells.Item(140, 1).Formula = "'2026-01-28"
$ws.Cells.Item(140, 2).Value = "17:41:49"
$ws.Cells.Item(140, 3).Value = "17:00"
$ws.Cells.Item(140, 4).Value = "Bedroom"
$ws.Cells.Item(140, 5).Value = "In Bed | HR=56 | BR=8"
$ws.Cells.Item(140, 6).Value = "Occupied"

$ws.Cells.Item(141, 1).Formula = "'2026-01-28"
$ws.Cells.Item(141, 2).Value = "17:41:50"
$ws.Cells.Item(141, 3).Value = "17:00"
$ws.Cells.Item(141, 4).Value = "Bedroom"
$ws.Cells.Item(141, 5).Value = "In Bed | HR=58 | BR=10"
$ws.Cells.Item(141, 6).Value = "Occupied"

$ws.Cells.Item(142, 1).Formula = "'2026-01-28"
$ws.Cells.Item(142, 2).Value = "17:41:52"
$ws.Cells.Item(142, 3).Value = "17:00"
$ws.Cells.Item(142, 4).Value = "Bedroom"
$ws.Cells.Item(142, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(142, 6).Value = "Occupied"

$ws.Cells.Item(143, 1).Formula = "'2026-01-28"
$ws.Cells.Item(143, 2).Value = "17:41:52"
$ws.Cells.Item(143, 3).Value = "17:00"
$ws.Cells.Item(143, 4).Value = "Bedroom"
$ws.Cells.Item(143, 5).Value = "In Bed | HR=55 | BR=7"
$ws.Cells.Item(143, 6).Value = "Occupied"

$ws.Cells.Item(144, 1).Formula = "'2026-01-28"
$ws.Cells.Item(144, 2).Value = "17:41:53"
$ws.Cells.Item(144, 3).Value = "17:00"
$ws.Cells.Item(144, 4).Value = "Bedroom"
$ws.Cells.Item(144, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(144, 6).Value = "Occupied"

$ws.Cells.Item(145, 1).Formula = "'2026-01-28"
$ws.Cells.Item(145, 2).Value = "17:41:55"
$ws.Cells.Item(145, 3).Value = "17:00"
$ws.Cells.Item(145, 4).Value = "Bedroom"
$ws.Cells.Item(145, 5).Value = "In Bed | HR=65 | BR=17"
$ws.Cells.Item(145, 6).Value = "Occupied"

$ws.Cells.Item(146, 1).Formula = "'2026-01-28"
$ws.Cells.Item(146, 2).Value = "17:41:55"
$ws.Cells.Item(146, 3).Value = "17:00"
$ws.Cells.Item(146, 4).Value = "Bedroom"
$ws.Cells.Item(146, 5).Value = "In Bed | HR=91 | BR=43"
$ws.Cells.Item(146, 6).Value = "Occupied"

$ws.Cells.Item(147, 1).Formula = "'2026-01-28"
$ws.Cells.Item(147, 2).Value = "17:41:56"
$ws.Cells.Item(147, 3).Value = "17:00"
$ws.Cells.Item(147, 4).Value = "Bedroom"
$ws.Cells.Item(147, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(147, 6).Value = "Occupied"

$ws.Cells.Item(148, 1).Formula = "'2026-01-28"
$ws.Cells.Item(148, 2).Value = "17:41:58"
$ws.Cells.Item(148, 3).Value = "17:00"
$ws.Cells.Item(148, 4).Value = "Bedroom"
$ws.Cells.Item(148, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(148, 6).Value = "Occupied"

# --- Camera sheet: add row 5 ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(5, 1).Formula = "'2026-01-28"
$ws.Cells.Item(5, 2).Value = "17:41:17"
$ws.Cells.Item(5, 3).Value = "17:00"
$ws.Cells.Item(5, 4).Value = "Living Room Main Door"
$ws.Cells.Item(5, 5).Value = "Image Captured"
$ws.Cells.Item(5, 6).Value = "Active"

